$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: role/kind metadata per column.
# Column A (tipo-de-estudios-realizados) is now curated as a measure (was a dimension).
$ws.Range("A2").Value = "iaest-measure:tipo-de-estudios-realizados"
# Column B (municipio-nombre) is now curated as an SDMX dimension (was a measure).
$ws.Range("B2").Value = "sdmx-dimension:refArea"

# Row 3: dim/medida swap between columns A and B.
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"

# Row 4: datatype swap between columns A and B.
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "URI-Municipio"

# Row 5 no longer applies (the mapping-file reference belonged to the old
# dimension column A, which is now a measure) - remove it entirely.
$ws.Range("A5").Delete()
